$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell B1 ("Light power") gets its text extended.
$ws.Range("B1").Value = "Light power22222222222222"

# Row 1 wraps to two lines now that the header text is longer; adjust its height.
$ws.Rows.Item(1).RowHeight = 58

# Move the active selection to B1 (was B9)
$ws.Range("B1").Select()
